# Update column F (dSF) values per repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = -6
    3  = -1
    4  = -2
    6  = -1
    7  = 1
    8  = 5
    9  = 2
    10 = 3
    11 = -2
    12 = 1
    13 = -3
    16 = -2
    17 = -1
    18 = -2
    19 = 1
    20 = -1
    21 = 3
    22 = 5
    23 = -1
    25 = -3
    26 = -2
    27 = -4
    28 = 5
    29 = -4
    30 = -1
    31 = -2
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
